$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("G2").Value = 2.7
$ws.Range("H2").Value = 2.9
$ws.Range("I2").Value = 2.9
$ws.Range("K2").Value = 1.91
$ws.Range("L2").Value = 3.75
$ws.Range("O2").Value = 1.53
$ws.Range("P2").Value = 2.38
$ws.Range("Q2").Value = 2.7
$ws.Range("R2").Value = 1.44
$ws.Range("U2").Value = 2.1
$ws.Range("V2").Value = 1.67
$ws.Range("W2").Value = 6.5
$ws.Range("X2").Value = 11
$ws.Range("Z2").Value = 26
$ws.Range("AE2").Value = 19
$ws.Range("AJ2").Value = 12
$ws.Range("AL2").Value = 29
$ws.Range("AR2").Value = 101
$ws.Range("AW2").Value = 4.75
$ws.Range("AY2").Value = 34
$ws.Range("BD2").Value = 151

# Row 6 updates
$ws.Range("I6").Value = 1.45
$ws.Range("L6").Value = 2.05
$ws.Range("AJ6").Value = 9
$ws.Range("AM6").Value = 29
$ws.Range("AR6").Value = 151
$ws.Range("AX6").Value = 7.5
$ws.Range("AZ6").Value = 23

# Row 7 updates
$ws.Range("H7").Value = 3.1
$ws.Range("I7").Value = 3.5
$ws.Range("J7").Value = 3
$ws.Range("K7").Value = 2
$ws.Range("S7").Value = 1.5
$ws.Range("T7").Value = 2.5
$ws.Range("W7").Value = 6.5
$ws.Range("AA7").Value = 21
$ws.Range("AC7").Value = 7.5
$ws.Range("AD7").Value = 6
$ws.Range("AE7").Value = 15
$ws.Range("AH7").Value = 9
$ws.Range("AI7").Value = 17
$ws.Range("AJ7").Value = 13
$ws.Range("AK7").Value = 41
$ws.Range("AP7").Value = 26
$ws.Range("AT7").Value = 2.5
$ws.Range("AV7").Value = 67
$ws.Range("BA7").Value = 101
